$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TAYDA ORDER")

# ---------------------------------------------------------------
# 1) Insert a brand-new BOM line (UV printing service) right below
#    the DRILL1 row, i.e. as the new row 5, pushing everything else
#    down by one row.
# ---------------------------------------------------------------
$ws.Rows("5:5").Insert()
$ws.Rows("5:5").RowHeight = 14.25

# copy the same number formatting / style used by the surrounding
# data rows (row 6 is the row that used to be row 5 before the
# insert, and is a normal -- non-total -- data row)
$ws.Range("A6:F6").Copy()
$ws.Range("A5:F5").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# fill in the new row's content (order matches the order the new
# strings were appended to the workbook: Product, RefDes, Part No)
$ws.Range("B5").Value = "125B ENCLOSURE FACE UV PRINTING SERVICE SKU: A-5165-CST-UV1"
$ws.Range("A5").Value = "UV1"
$ws.Range("C5").Value = "TAYDA UV PRINTING TEMPLATE"
$ws.Range("D5").Formula = '=(LEN(A5)-LEN(SUBSTITUTE(A5,",","")) + 1)'
$ws.Range("E5").Value = 4
$ws.Range("F5").Formula = "='TAYDA ORDER'!`$E5*'TAYDA ORDER'!`$D5"

# ---------------------------------------------------------------
# 2) Fix the text on the existing DRILL1 row (row 4): the template
#    name is shortened from "TAYDA DRILL TEMPLATE 1590BB" to
#    "TAYDA DRILL TEMPLATE".
# ---------------------------------------------------------------
$ws.Range("C4").Value = "TAYDA DRILL TEMPLATE"

# ---------------------------------------------------------------
# 3) Grow the Table_3 ListObject so the new row becomes part of the
#    table (and the totals row formulas / range shift down too).
# ---------------------------------------------------------------
$tbl = $ws.ListObjects.Item("Table_3")
$tbl.Resize($ws.Range("A1:F15"))

# ---------------------------------------------------------------
# 4) Cosmetic selection change recorded in the workbook XML.
# ---------------------------------------------------------------
$ws.Range("B10").Select()

Write-Host "done"
